# LCV duzeltmeleri yapildi, HDV ilk model eklendi
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "HCV_MTV" sheet at the very end of the workbook (after
#    "lcv_co2_aralikli_mtv"), copying header formatting from LCV_MTV_oranlari
# ---------------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("LCV_MTV_oranlari")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "HCV_MTV"

# Copy header (row 1) cell formatting (style, fill, borders, alignment) from
# the LCV_MTV_oranlari sheet, whose header uses the same look.
$srcSheet.Range("A1:AJ1").Copy()
$new.Range("A1:AJ1").PasteSpecial(-4122)

# Header captions
$header = @("mtv_grubu","govde_tipi","agirlik_min","agirlik_max","oturma_yeri_min","oturma_yeri_max","sene_1","sene_2","sene_3","sene_4","sene_5","sene_6","sene_7","sene_8","sene_9","sene_10","sene_11","sene_12","sene_13","sene_14","sene_15","sene_16","sene_17","sene_18","sene_19","sene_20","sene_21","sene_22","sene_23","sene_24","sene_25","sene_26","sene_27","sene_28","sene_29","sene_30")
for ($i = 0; $i -lt $header.Length; $i++) {
    $new.Cells.Item(1, $i + 1).Value = $header[$i]
}
$new.Rows.Item(1).RowHeight = 47.25

# Data rows (A2:AJ11)
$data = @(
    @(1,"kamyon",0,1500,0,99999,1033,1033,1033,1033,1033,1033,686,686,686,686,686,686,686,686,686,337,337,337,337,337,337,337,337,337,337,337,337,337,337,337),
    @(2,"kamyon",1501,3500,0,99999,2091,2091,2091,2091,2091,2091,1212,1212,1212,1212,1212,1212,1212,1212,1212,686,686,686,686,686,686,686,686,686,686,686,686,686,686,686),
    @(3,"kamyon",3501,5000,0,99999,3141,3141,3141,3141,3141,3141,2614,2614,2614,2614,2614,2614,2614,2614,2614,1033,1033,1033,1033,1033,1033,1033,1033,1033,1033,1033,1033,1033,1033,1033),
    @(4,"kamyon",5001,10000,0,99999,3483,3483,3483,3483,3483,3483,2962,2962,2962,2962,2962,2962,2962,2962,2962,1389,1389,1389,1389,1389,1389,1389,1389,1389,1389,1389,1389,1389,1389,1389),
    @(5,"kamyon",10001,20000,0,99999,4191,4191,4191,4191,4191,4191,3487,3487,3487,3487,3487,3487,3487,3487,3487,2091,2091,2091,2091,2091,2091,2091,2091,2091,2091,2091,2091,2091,2091,2091),
    @(6,"kamyon",20001,999999,0,99999,5242,5242,5242,5242,5242,5242,4191,4191,4191,4191,4191,4191,4191,4191,4191,2436,2436,2436,2436,2436,2436,2436,2436,2436,2436,2436,2436,2436,2436,2436),
    @(7,"otobus",0,999999,0,25,727,727,727,727,727,727,434,434,434,434,434,434,434,434,434,190,190,190,190,190,190,190,190,190,190,190,190,190,190,190),
    @(8,"otobus",0,999999,26,35,871,871,871,871,871,871,727,727,727,727,727,727,727,727,727,288,288,288,288,288,288,288,288,288,288,288,288,288,288,288),
    @(9,"otobus",0,999999,36,45,970,970,970,970,970,970,823,823,823,823,823,823,823,823,823,384,384,384,384,384,384,384,384,384,384,384,384,384,384,384),
    @(10,"otobus",0,999999,46,99999,1163,1163,1163,1163,1163,1163,970,970,970,970,970,970,970,970,970,581,581,581,581,581,581,581,581,581,581,581,581,581,581,581)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $new.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}

# Column widths (engine stores OOXML width = ColumnWidth + 5/6, so back out
# the offset to land on the target stored widths of 12.5 / 13)
$new.Columns.Item(5).ColumnWidth = 11.666666666666666
$new.Columns.Item(6).ColumnWidth = 12.166666666666666

# View/selection state for the new sheet (becomes the active/selected tab)
$new.Range("U22").Select()

# ---------------------------------------------------------------------------
# 2) LCV_MTV_oranlari: selection now spans the header row A1:AJ1
# ---------------------------------------------------------------------------
$srcSheet.Range("A1:AJ1").Select()

# ---------------------------------------------------------------------------
# 3) lcv_co2_aralikli_mtv: selection back to D4 (tabSelected flag moves to
#    the new HCV_MTV sheet automatically once it is added/activated)
# ---------------------------------------------------------------------------
$oldLast = $wb.Worksheets.Item("lcv_co2_aralikli_mtv")
$oldLast.Range("D4").Select()

# Re-activate the new sheet so it ends up as the active/selected tab
$new.Activate()
$new.Range("U22").Select()
